$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 920.08826
$ws.Range("I6").Value = 35.5
$ws.Range("J6").Value = 975.375
$ws.Range("K6").Value = 106.5
$ws.Range("L6").Value = 2926.125
$ws.Range("M6").Value = 5.5
$ws.Range("N6").Value = -3150.125
$ws.Range("H31").Value = 2699.6
$ws.Range("I31").Value = 624.75
$ws.Range("K31").Value = 1874.25
$ws.Range("M31").Value = -1644.25
$ws.Range("H86").Value = 1002726
$ws.Range("J86").Value = 3002
$ws.Range("L86").Value = 3002
$ws.Range("N86").Value = -5248
$ws.Range("H89").Value = 1002726
$ws.Range("J89").Value = 3002
$ws.Range("L89").Value = 15010
$ws.Range("N89").Value = -26242
$ws.Range("H116").Value = 3549.1304
$ws.Range("I116").Value = 3443.4443
$ws.Range("K116").Value = 3443.4443
$ws.Range("M116").Value = -1.444300000000112
$ws.Range("H140").Value = 89996.5
$ws.Range("J140").Value = 89996.5
$ws.Range("L140").Value = 89996.5
$ws.Range("N140").Value = -100356.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 2000
$ws.Range("J3").Value = 2000
$ws.Range("L3").Value = 2000
$ws.Range("N3").Value = -2230
$ws.Range("H6").Value = 5204200.5
$ws.Range("J6").Value = 6503500
$ws.Range("L6").Value = 6503500
$ws.Range("N6").Value = -6503846
$ws.Range("H11").Value = 20059800
$ws.Range("I11").Value = 25049500
$ws.Range("K11").Value = 25049500
$ws.Range("M11").Value = -25049356
$ws.Range("H32").Value = 7884.68
$ws.Range("I32").Value = 664.36584
$ws.Range("K32").Value = 664.36584
$ws.Range("M32").Value = -377.36584

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 7040.8
$ws.Range("I11").Value = 1100
$ws.Range("J11").Value = 11001.333
$ws.Range("K11").Value = 1100
$ws.Range("L11").Value = 11001.333
$ws.Range("M11").Value = -960
$ws.Range("N11").Value = -11281.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = ""
$ws.Range("H3").Value = 1333335.4
$ws.Range("I3").Value = 1000003
$ws.Range("K3").Value = 1000003
$ws.Range("M3").Value = -999890
$ws.Range("H4").Value = 450638
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 600684
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 600684
$ws.Range("M4").Value = -388
$ws.Range("N4").Value = -600908
$ws.Range("H5").Value = 16011
$ws.Range("J5").Value = 16011
$ws.Range("L5").Value = 16011
$ws.Range("N5").Value = -16235
$ws.Range("H8").Value = 750
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 750
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 750
$ws.Range("M8").Value = ""
$ws.Range("N8").Value = -1030
$ws.Range("H10").Value = 2500.5715
$ws.Range("I10").Value = 265.33334
$ws.Range("J10").Value = 4177
$ws.Range("K10").Value = 265.33334
$ws.Range("L10").Value = 4177
$ws.Range("M10").Value = -126.33334
$ws.Range("N10").Value = -4455
$ws.Range("H11").Value = 5000
$ws.Range("J11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("N11").Value = -5280
$ws.Range("H13").Value = 2000
$ws.Range("J13").Value = 2000
$ws.Range("L13").Value = 2000
$ws.Range("N13").Value = -2278
$ws.Range("H14").Value = 20840.666
$ws.Range("J14").Value = 20840.666
$ws.Range("L14").Value = 20840.666
$ws.Range("N14").Value = -21180.666
$ws.Range("H15").Value = 100
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").Value = ""
$ws.Range("H19").Value = 386.81818
$ws.Range("I19").Value = 110
$ws.Range("K19").Value = 110
$ws.Range("M19").Value = 60
$ws.Range("H24").Value = 386.81818
$ws.Range("I24").Value = 110
$ws.Range("K24").Value = 110
$ws.Range("M24").Value = 60
$ws.Range("H25").Value = 1262.7778
$ws.Range("I25").Value = 1043.3334
$ws.Range("J25").Value = 1701.6666
$ws.Range("K25").Value = 1043.3334
$ws.Range("L25").Value = 1701.6666
$ws.Range("M25").Value = -869.3334
$ws.Range("N25").Value = -2049.6666
$ws.Range("H31").Value = 5048.3716
$ws.Range("J31").Value = 4382.353
$ws.Range("L31").Value = 4382.353
$ws.Range("N31").Value = -4972.353
$ws.Range("H34").Value = 5048.3716
$ws.Range("J34").Value = 4382.353
$ws.Range("L34").Value = 4382.353
$ws.Range("N34").Value = -4786.353

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 6598.3335
$ws.Range("I14").Value = 6598.3335
$ws.Range("K14").Value = 19795.0005
$ws.Range("M14").Value = -19622.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 25001.5
$ws.Range("J10").Value = 25001.5
$ws.Range("L10").Value = 25001.5
$ws.Range("N10").Value = -25339.5
$ws.Range("H80").Value = 4919.1816
$ws.Range("I80").Value = 3053.25
$ws.Range("J80").Value = 5985.4287
$ws.Range("K80").Value = 3053.25
$ws.Range("L80").Value = 5985.4287
$ws.Range("M80").Value = -2055.25
$ws.Range("N80").Value = -7981.4287
$ws.Range("H83").Value = 4919.1816
$ws.Range("I83").Value = 3053.25
$ws.Range("J83").Value = 5985.4287
$ws.Range("K83").Value = 15266.25
$ws.Range("L83").Value = 29927.1435
$ws.Range("M83").Value = -10274.25
$ws.Range("N83").Value = -39911.14350000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 500875
$ws.Range("J5").Value = 500875
$ws.Range("L5").Value = 500875
$ws.Range("N5").Value = -501099
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = ""
$ws.Range("H13").Value = 400
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""
